$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix header typo: "Cantidad de clases" -> "Lecturas"
$ws.Range("C1").Value = "Lecturas"

# Column C (rows 2-8) previously held numeric class counts; replace with "."
# to match the other placeholder columns (D, E)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = "."
}

# Update view: zoom to 150% and move the active selection to E2
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("E2").Select()
